$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Row 31 ---
$lo.ListRows.Add() | Out-Null
$ws.Range("A30:F30").Copy()
$ws.Range("A31:F31").Insert(-4121)
$ws.Cells.Item(31,1).Value = 43356
$ws.Cells.Item(31,2).Value = 0.84375
$ws.Cells.Item(31,3).Value = 0.89097222222222217
$ws.Cells.Item(31,4).Formula = "=(C31-B31)* 1440"
$ws.Cells.Item(31,5).Formula = "=IF(C31>B31, (C31-B31)*1440, (B31-C31)*1440)"
$ws.Cells.Item(31,6).Formula = "=ABS((C31-B31)*1440)"

# --- Row 32 ---
$lo.ListRows.Add() | Out-Null
$ws.Range("A31:F31").Copy()
$ws.Range("A32:F32").Insert(-4121)
$ws.Cells.Item(32,1).Value = 43356
$ws.Cells.Item(32,2).Value = 0.98402777777777783
$ws.Cells.Item(32,3).Value = 0.99930555555555556
$ws.Cells.Item(32,4).Formula = "=(C32-B32)* 1440"
$ws.Cells.Item(32,5).Formula = "=IF(C32>B32, (C32-B32)*1440, (B32-C32)*1440)"
$ws.Cells.Item(32,6).Formula = "=ABS((C32-B32)*1440)"

# --- Row 33 ---
$lo.ListRows.Add() | Out-Null
$ws.Range("A32:F32").Copy()
$ws.Range("A33:F33").Insert(-4121)
$ws.Cells.Item(33,1).Value = 43357
$ws.Cells.Item(33,2).Value = 0
$ws.Cells.Item(33,3).Value = 0.27777777777777779
$ws.Cells.Item(33,4).Formula = "=(C33-B33)* 1440"
$ws.Cells.Item(33,5).Formula = "=IF(C33>B33, (C33-B33)*1440, (B33-C33)*1440)"
$ws.Cells.Item(33,6).Formula = "=ABS((C33-B33)*1440)"

$excel.CutCopyMode = $false
$ws.Range("D33").Select()
